# Commit: "Fruta / hortaliza, semanal"
# Inserts one new daily price record for "Zapallo italiano" (Terminal La
# Palmera de La Serena) right after the existing row 100 (Excel row 101),
# pushing every subsequent record down by one row. The sheet's used range
# grows from A1:R194 to A1:R195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 101:194 down to 102:195, leaving a blank row 101 that
# inherits the formatting (incl. the date number-format on column D)
# of the row above it, matching Excel's native Insert behaviour.
$ws.Rows.Item(101).Insert()

# Populate the newly-opened row 101 with the new record.
$ws.Range("A101").Value = 8
$ws.Range("B101").Value = "Terminal La Palmera de La Serena"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44484
$ws.Range("E101").Value = 4
$ws.Range("F101").Value = 100112032
$ws.Range("G101").Value = "Zapallo italiano"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 500
$ws.Range("K101").Value = 13000
$ws.Range("L101").Value = 14000
$ws.Range("M101").Value = 13500
$ws.Range("N101").Value = "$/caja 70 unidades"
$ws.Range("O101").Value = "Provincia de Limarí"
$ws.Range("P101").Value = 193
$ws.Range("Q101").Value = 70
$ws.Range("R101").Value = "Hortaliza"
